$d = $word.ActiveDocument

# Document layout (before edit), one paragraph each:
#   1: "⚡️🚀המאמר היומי של מייק 27.07.24: ⚡️🚀"
#   2: paper title
#   3: first body paragraph (kept, text replaced)
#   4: second body paragraph (to be removed entirely)
#   5: third body paragraph (to be removed entirely)
#   6: arxiv URL (kept, text replaced)
#
# Use direct Range.Text assignment (not Find/Replace) so Word's
# AutoCorrect "smart quotes" feature does not mangle the straight
# apostrophe in the new body text.

# 1) Update the date in the daily title line: 27.07.24 -> 26.07.24
$d.Paragraphs.Item(1).Range.Text = "⚡️🚀המאמר היומי של מייק 26.07.24: ⚡️🚀"

# 2) Replace the paper title line with the new paper title
$d.Paragraphs.Item(2).Range.Text = "Questionable practices in machine learning"

# 3) Replace the first body paragraph with the new (short) review text
$newBody = "הסקירה היום תהיה ממש קצרה. המאמר המסוקר דן בפרקטיקות פסולות שעלולות להכשיל אתכם במהלך פיתוח של המודלים שלכם. רוב הפרקטיקות הרעות שנזכרו במאמר נראות לחוקרי ML מנוסים די טריוויאליות ודי ברור למה לא כדאי להשתמש בהן. בין אלו ניתן למנות אימון על טסט סט, בחירה של בייסליין חלש להשוואה, הסקת מסקנות על אימון אחד בלבד של המודל, אימון על דאטה דומה מאוד לבנצ'מארק וכדומה. אבל ניתן למצוא גם דברים פחות טריוויאליים שחלקם לא ידעתי. "
$d.Paragraphs.Item(3).Range.Text = $newBody

# 4) Delete the two paragraphs that explained BPE tokenization / the linear
#    algorithm in detail (paragraphs 4 and 5, just before the URL
#    paragraph). Deleting index 4 twice removes both, since later
#    paragraphs shift up by one each time a paragraph is removed.
$d.Paragraphs.Item(4).Range.Delete() | Out-Null
$d.Paragraphs.Item(4).Range.Delete() | Out-Null

# 5) Update the arxiv link (now paragraph 4 after the deletions above)
$d.Paragraphs.Item(4).Range.Text = "https://www.arxiv.org/abs/2407.12220"
